# Write obj file address to excel file
# Adds three new rows (14-16) to the "Submenu" sheet describing the
# "house" object model and its .obj file paths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Submenu")

# Row 14: house / SkyBlue / 15 / C:/Users/xzjzb/PycharmProjects/MapEditor/data/model/house/house.obj
$ws.Range("A14").Value = "house"
$ws.Range("B14").Value = "SkyBlue"
$ws.Range("C14").Value = 15
$ws.Range("D14").Value = "C:/Users/xzjzb/PycharmProjects/MapEditor/data/model/house/house.obj"

# Row 15: house2 / Maroon / 11 / /data/model/house/house.obj
$ws.Range("A15").Value = "house2"
$ws.Range("B15").Value = "Maroon"
$ws.Range("C15").Value = 11
$ws.Range("D15").Value = "/data/model/house/house.obj"

# Row 16: 建物house / Tan / 16 / data/model/house/house.obj
$ws.Range("A16").Value = "建物house"
$ws.Range("B16").Value = "Tan"
$ws.Range("C16").Value = 16
$ws.Range("D16").Value = "data/model/house/house.obj"
